$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.918.85"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "2.265.95"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "301.93"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").Value = "92.08"
$ws.Range("E6").Value = "  +5.76%  "
$ws.Range("E7").Value = "  +3.32%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "0.485"
$ws.Range("E9").Value = "  +4.25%  "
$ws.Range("D10").Value = "54.44"
$ws.Range("E10").Value = "  +8.06%  "
$ws.Range("D11").Value = "32.29"
$ws.Range("E11").Value = "  +6.03%  "
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "2.618.41"
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").Value = "14.15"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "2.270.37"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "0.759"
$ws.Range("E18").Value = "  +3.63%  "
$ws.Range("D19").Value = "41.842.77"
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("E20").Value = "  +8.93%  "
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("E22").Value = "  +3.42%  "
$ws.Range("D23").Value = "67.02"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Value = "241.86"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +3.58%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.90"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").Value = "23.94"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("D31").Value = "159.45"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "33.83"
$ws.Range("E32").Value = "  +6.30%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("D35").Value = "0.0746"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("D36").Value = "3.08"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.60"
$ws.Range("E38").Value = "  +9.10%  "
$ws.Range("E39").Value = "  +4.98%  "
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D42").Value = "3.92"
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("D43").Value = "2.067.88"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "19.65"
$ws.Range("E44").Value = "  +9.24%  "
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("D46").Value = "10.22"
$ws.Range("E46").Value = "  +4.46%  "
$ws.Range("E47").Value = "  +7.49%  "
$ws.Range("D48").Value = "2.03"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("E49").Value = "  +3.74%  "
$ws.Range("D50").Value = "1.14"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("D51").Value = "51.76"
$ws.Range("E51").Value = "  +5.51%  "
